$d = $word.ActiveDocument

# 1. Fix the typo "write prints" -> "writeprints" (the paragraph about
#    authorship verification talks about comparing "writeprints").
$d.Content.Find.Execute("write prints", $true, $false, $false, $false, $false, $true, 1, $false, "writeprints", 2)

# 2. Everything that used to follow the sentence "... to verify authorship."
#    (the link-detection / NLP / Stanford-parser notes, running all the way to
#    the end of the document body) gets removed.
$anchor = $d.Content
$anchor.Find.Execute("to verify authorship.")
$anchor.MoveEnd(1, 1)          # wdCharacter - also swallow the paragraph mark
$cutStart = $anchor.End

$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$cutEnd = $lastParagraph.Range.End

$deleteRange = $d.Range($cutStart, $cutEnd)
$deleteRange.Delete()
